# Update order #14 (row 14) on the "All Orders" sheet:
#   - Status (H14) changes from NEW to CANCELLED
#   - Cancel Reason (M14) is set to "test order"
$wb = $excel.ActiveWorkbook

$ordersSheet = $wb.Worksheets.Item("All Orders")
$ordersSheet.Range("H14").Value = "CANCELLED"
$ordersSheet.Range("M14").Value = "test order"

# Update the "Daily Summary" sheet totals for 2026-01-13 (row 4) to reflect
# the newly cancelled order:
#   - Cancelled (D4): 8 -> 9
#   - Revenue (E4):   95 -> 80
#   - Pending (G4):   95 -> 80
$summarySheet = $wb.Worksheets.Item("Daily Summary")
$summarySheet.Range("D4").Value = 9
$summarySheet.Range("E4").Value = 80
$summarySheet.Range("G4").Value = 80
